$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.715275406837463
$ws.Range("B1").Value = 2.761646509170532
$ws.Range("C1").Value = 2.994082450866699
$ws.Range("D1").Value = 3.375754356384277
$ws.Range("E1").Value = 2.062291145324707
